$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column M1 header "PAID"
$ws.Range("M1").Value = "PAID"

# Force text storage (no numeric auto-inference) for the numeric-looking
# Donation Amount (F), Phone (H) and Zip (L) columns across the new rows,
# matching the existing t="str" cells used throughout the sheet.
$ws.Range("F15:F24").NumberFormat = "@"
$ws.Range("H15:H24").NumberFormat = "@"
$ws.Range("L15:L24").NumberFormat = "@"

# Row 15
$ws.Range("A15").Value = "REC-1741969487529-984"
$ws.Range("B15").Value = 14
$ws.Range("C15").Value = "2025-03-14T16:24:47.546Z"
$ws.Range("D15").Value = "Deepak"
$ws.Range("E15").Value = "Adhikari"
$ws.Range("F15").Value = "988"
$ws.Range("G15").Value = "sjahbfkjadhkjfa@gmail.com"
$ws.Range("H15").Value = "3477712375"
$ws.Range("I15").Value = "11 alpine ln"
$ws.Range("J15").Value = "Hicksville"
$ws.Range("K15").Value = "NY"
$ws.Range("L15").Value = "11801"

# Row 16
$ws.Range("A16").Value = "REC-1741969570079-329"
$ws.Range("B16").Value = 15
$ws.Range("C16").Value = "2025-03-14T16:26:10.085Z"
$ws.Range("D16").Value = "Deepak"
$ws.Range("E16").Value = "Adhikari"
$ws.Range("F16").Value = "988"
$ws.Range("G16").Value = "sjahbfkjadhkjfa@gmail.com"
$ws.Range("H16").Value = "3477712375"
$ws.Range("I16").Value = "11 alpine ln"
$ws.Range("J16").Value = "Hicksville"
$ws.Range("K16").Value = "NY"
$ws.Range("L16").Value = "11801"

# Row 17
$ws.Range("A17").Value = "REC-1741969650869-422"
$ws.Range("B17").Value = 16
$ws.Range("C17").Value = "2025-03-14T16:27:30.874Z"
$ws.Range("D17").Value = "Deepak"
$ws.Range("E17").Value = "Adhikari"
$ws.Range("F17").Value = "988"
$ws.Range("G17").Value = "sjahbfkjadhkjfa@gmail.com"
$ws.Range("H17").Value = "3477712375"
$ws.Range("I17").Value = "11 alpine ln"
$ws.Range("J17").Value = "Hicksville"
$ws.Range("K17").Value = "NY"
$ws.Range("L17").Value = "11801"

# Row 18
$ws.Range("A18").Value = "REC-1741969739843-247"
$ws.Range("B18").Value = 17
$ws.Range("C18").Value = "2025-03-14T16:28:59.849Z"
$ws.Range("D18").Value = "Deepak"
$ws.Range("E18").Value = "Adhikari"
$ws.Range("F18").Value = "333333"
$ws.Range("G18").Value = "dadhikari856@gmail.com"
$ws.Range("H18").Value = "3477712375"
$ws.Range("I18").Value = "11 alpine ln"
$ws.Range("J18").Value = "Hicksville"
$ws.Range("K18").Value = "NY"
$ws.Range("L18").Value = "11801"

# Row 19
$ws.Range("A19").Value = "REC-1741970623658-586"
$ws.Range("B19").Value = 18
$ws.Range("C19").Value = "2025-03-14T16:43:43.665Z"
$ws.Range("D19").Value = "Deepak"
$ws.Range("E19").Value = "Adhikari"
$ws.Range("F19").Value = "2222333"
$ws.Range("G19").Value = "dadhikari856@gmail.com"
$ws.Range("H19").Value = "3477712375"
$ws.Range("I19").Value = "11 alpine ln"
$ws.Range("J19").Value = "Hicksville"
$ws.Range("K19").Value = "NY"
$ws.Range("L19").Value = "11801"

# Row 20
$ws.Range("A20").Value = "REC-1741970747227-173"
$ws.Range("B20").Value = 19
$ws.Range("C20").Value = "2025-03-14T16:45:47.232Z"
$ws.Range("D20").Value = "Deepak"
$ws.Range("E20").Value = "Adhikari"
$ws.Range("F20").Value = "2222333"
$ws.Range("G20").Value = "dadhikari856@gmail.com"
$ws.Range("H20").Value = "3477712375"
$ws.Range("I20").Value = "11 alpine ln"
$ws.Range("J20").Value = "Hicksville"
$ws.Range("K20").Value = "NY"
$ws.Range("L20").Value = "11801"

# Row 21
$ws.Range("A21").Value = "REC-1741971092861-971"
$ws.Range("B21").Value = 20
$ws.Range("C21").Value = "2025-03-14T16:51:32.867Z"
$ws.Range("D21").Value = "Deepak"
$ws.Range("E21").Value = "Adhikari"
$ws.Range("F21").Value = "2222333"
$ws.Range("G21").Value = "dadhikari856@gmail.com"
$ws.Range("H21").Value = "3477712375"
$ws.Range("I21").Value = "11 alpine ln"
$ws.Range("J21").Value = "Hicksville"
$ws.Range("K21").Value = "NY"
$ws.Range("L21").Value = "11801"

# Row 22
$ws.Range("A22").Value = "REC-1741971120418-780"
$ws.Range("B22").Value = 21
$ws.Range("C22").Value = "2025-03-14T16:52:00.420Z"
$ws.Range("D22").Value = "Deepak"
$ws.Range("E22").Value = "Adhikari"
$ws.Range("F22").Value = "6666666666"
$ws.Range("G22").Value = "dadhikari856@gmail.com"
$ws.Range("H22").Value = "3477712375"
$ws.Range("I22").Value = "11 alpine ln"
$ws.Range("J22").Value = "Hicksville"
$ws.Range("K22").Value = "NY"
$ws.Range("L22").Value = "11801"

# Row 23
$ws.Range("A23").Value = "REC-1741971739730-822"
$ws.Range("B23").Value = 22
$ws.Range("C23").Value = "2025-03-14T17:02:19.735Z"
$ws.Range("D23").Value = "Deepak"
$ws.Range("E23").Value = "Adhikari"
$ws.Range("F23").Value = "132"
$ws.Range("G23").Value = "deepak-adhikari@hotmail.com"
$ws.Range("H23").Value = "8567768105"
$ws.Range("I23").Value = "11 alpine ln"
$ws.Range("J23").Value = "Hicksville"
$ws.Range("K23").Value = "Alabama"
$ws.Range("L23").Value = "11801"

# Row 24
$ws.Range("A24").Value = "REC-1741971785782-752"
$ws.Range("B24").Value = 23
$ws.Range("C24").Value = "2025-03-14T17:03:05.790Z"
$ws.Range("D24").Value = "Deepak"
$ws.Range("E24").Value = "Adhikari"
$ws.Range("F24").Value = "132"
$ws.Range("G24").Value = "deepak-adhikari@taptap.com"
$ws.Range("H24").Value = "8567768105"
$ws.Range("I24").Value = "11 alpine ln"
$ws.Range("J24").Value = "Hicksville"
$ws.Range("K24").Value = "Alabama"
$ws.Range("L24").Value = "11801"

